$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Unprotect the sheet to allow edits (it was protected in the original file)
$ws.Unprotect("lido")

$ws.Range("A10").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-05-17 for illustrative purposes only and are subject to change."

$ws.Range("D2").Value = 0.4754863459309436
$ws.Range("E2").Value = -0.0003929273084479767

$ws.Range("D3").Value = 0.3432931768373103
$ws.Range("E3").Value = 0.0005585552038727126

$ws.Range("D4").Value = 0.09576347123224024
$ws.Range("E4").Value = 0.0008984725965859308

$ws.Range("D5").Value = 0.05359730667541825
$ws.Range("E5").Value = -0.001834651989450653

$ws.Range("D6").Value = 0.03185969932408764
$ws.Range("E6").Value = 0.01851851851851838

$ws.Range("D7").Value = 1
$ws.Range("E7").Value = 0.0005826195015286029

# Restore sheet protection to match the workbook's original protected state
$ws.Protect("lido", $true, $true, $true)
